$d = $word.ActiveDocument

# 1) Update the "Curso (semestre ideal)" line with the new semester numbers.
$d.Content.Find.Execute(
    "Curso (semestre ideal): EA (4), EB (3), EQD (3), EQN (4)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Curso (semestre ideal): EA (4), EB (5), EQD (4), EQN (5)", 2)

# 2) Split the single "Requisitos" bullet into two bullets:
#    replace "LOB1004 -  Cálculo II  (Requisito fraco)" with
#    "LOB1024 -  Mecânica  (Requisito fraco)" followed by a new run for
#    "LOB1052 -  Cálculo III  (Requisito fraco)", both ending in a line
#    break, inside the same paragraph.
$search = $d.Content
$found = $search.Find.Execute(
    "LOB1004 -  Cálculo II  (Requisito fraco)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Extend the found range by one character so it also covers the
# trailing manual line break (Chr(11)) that follows the run's text,
# without swallowing the paragraph mark itself.
$target = $d.Range($search.Start, $search.End + 1)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p>' +
       '<w:r><w:t>LOB1024 -  Mecânica  (Requisito fraco)</w:t><w:br/></w:r>' +
       '<w:r><w:t>LOB1052 -  Cálculo III  (Requisito fraco)</w:t><w:br/></w:r>' +
       '</w:p>' +
       '</w:body>' +
       '</w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)
